$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.323.46"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").Value = "2.046.16"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'228.82"
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'56.48"
$ws.Range("E8").Value = "  -3.34%  "

$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").Value = "'0.0785"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("D12").Value = "'14.70"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("D13").Value = "2.327.27"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").Value = "'20.61"
$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").Value = "'0.755"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "2.047.83"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "37.173.32"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").Value = "'6.05"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("D20").Value = "'69.28"
$ws.Range("E20").Value = "  -3.20%  "

$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -2.12%  "

$ws.Range("D22").Value = "'225.33"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'2.27"
$ws.Range("E25").Value = "  -5.13%  "

$ws.Range("D26").Value = "'9.66"
$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("D27").Value = "'166.31"
$ws.Range("E27").Value = "  -3.12%  "

$ws.Range("E28").Value = "  -7.25%  "

$ws.Range("D29").Value = "'18.97"
$ws.Range("E29").Value = "  -2.30%  "

$ws.Range("E30").Value = "  -4.14%  "

$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("D32").Value = "'4.52"
$ws.Range("E32").Value = "  -4.29%  "

$ws.Range("D33").Value = "'0.0613"
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("D34").Value = "'4.57"
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  -4.77%  "

$ws.Range("D39").Value = "'5.24"
$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("E40").Value = "  -4.78%  "

$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").Value = "1.478.20"
$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").Value = "'16.85"
$ws.Range("E43").Value = "  -0.88%  "

# Row 44: Cronos -> Aave (position swap with row 45)
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'96.20"
$ws.Range("E44").Value = "  -5.49%  "

# Row 45: Aave -> Cronos (position swap with row 44)
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0939"
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("E46").Value = "  +0.55%  "

# Row 47: FTXToken -> ARBITRUM (position swap with row 48)
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  -4.39%  "

# Row 48: ARBITRUM -> FTXToken (position swap with row 47)
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'3.94"
$ws.Range("E48").Value = "  -3.78%  "

$ws.Range("D49").Value = "'7.10"
$ws.Range("E49").Value = "  -3.68%  "

$ws.Range("D50").Value = "'2.92"
$ws.Range("E50").Value = "  -2.45%  "

$ws.Range("D51").Value = "2.229.86"
$ws.Range("E51").Value = "  -1.77%  "
